$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13.39237461141918
$ws.Range("E2").Value = 8.647496235152056
$ws.Range("F2").Value = 16.16266996806342
$ws.Range("G2").Value = 10.84599750390615
$ws.Range("H2").Value = 13.48906349216141
$ws.Range("I2").Value = 8.798273606621001
$ws.Range("J2").Value = 17.4833760584301
$ws.Range("K2").Value = 11.64918742987854
